$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# The USD Amount for the transaction in row 2 was corrected to 0.
$ws.Range("T2").Value = 0
